# Insert a new weekly price record as row 109 on the "Berenjena" sheet.
# This pushes the existing rows 109:207 down to 110:208 (dimension becomes A1:R208).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("109:109").Insert()

$ws.Cells.Item(109, 1).Value2 = 5
$ws.Cells.Item(109, 2).Value2 = "Macroferia Regional de Talca"
$ws.Cells.Item(109, 3).Value2 = "Maule"
$ws.Cells.Item(109, 4).Value2 = 45167
$ws.Cells.Item(109, 5).Value2 = 7
$ws.Cells.Item(109, 6).Value2 = 100112001
$ws.Cells.Item(109, 7).Value2 = "Berenjena"
$ws.Cells.Item(109, 8).Value2 = "Sin especificar"
$ws.Cells.Item(109, 9).Value2 = "Primera"
$ws.Cells.Item(109, 10).Value2 = 200
$ws.Cells.Item(109, 11).Value2 = 10000
$ws.Cells.Item(109, 12).Value2 = 10000
$ws.Cells.Item(109, 13).Value2 = 10000
$ws.Cells.Item(109, 14).Value2 = "`$/caja 50 unidades"
$ws.Cells.Item(109, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(109, 16).Value2 = 200
$ws.Cells.Item(109, 17).Value2 = 50
$ws.Cells.Item(109, 18).Value2 = "Hortaliza"
